# Update COVID-19 country stats and re-sort by "Casos totales" (Total cases)
# descending, then refresh the "last updated" timestamp in A1.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Data range boundaries (header is row 3, data starts row 4)
$firstDataRow = 4
$lastDataRow  = 219

# New statistics per country: Casos totales, Nuevos casos, Casos activos,
# Recuperados, Casos criticos, Muertes hoy, Muertes
$updates = [ordered]@{
    "Rusia"                              = @(871894, 5267, 676357, 180931, 0, 116, 14606)
    "Banglades"                          = @(249651, 2977, 143824, 102521, 0, 39, 3306)
    "Indonesia"                          = @(118753, 1882, 75645, 37587, 0, 69, 5521)
    "Filipinas"                          = @(119460, 3561, 66837, 50473, 0, 28, 2150)
    "Israel"                             = @(78514, 595, 53362, 24583, 0, 4, 569)
    "Singapur"                           = @(54555, 301, 47768, 6760, 0, 0, 27)
    "Polonia"                            = @(49515, 726, 35642, 12099, 0, 18, 1774)
    "Barein"                             = @(42514, 0, 39576, 2783, 0, 1, 155)
    "Austria"                            = @(21696, 130, 19596, 1381, 0, 0, 719)
    "Consejo Danes para los Refugiados"  = @(9309, 56, 8048, 1046, 0, 0, 215)
    "Finlandia"                          = @(7532, 20, 6980, 221, 0, 0, 331)
    "Hong Kong"                          = @(3850, 95, 2458, 1348, 0, 1, 44)
    "Eslovaquia"                         = @(2480, 63, 1824, 627, 0, 0, 29)
    "Lituania"                           = @(2171, 24, 1656, 434, 0, 0, 81)
    "Curazao"                            = @(31, 2, 28, 2, 0, 0, 1)
}

$searchRange = $ws.Range("A" + $firstDataRow + ":A" + $lastDataRow)

foreach ($country in $updates.Keys) {
    $hit = $searchRange.Find($country)
    if ($hit -eq $null) {
        Write-Output ("WARNING: country not found -> " + $country)
        continue
    }
    $r = $hit.Row
    $vals = $updates[$country]
    $ws.Range("B" + $r).Value = $vals[0]
    $ws.Range("C" + $r).Value = $vals[1]
    $ws.Range("D" + $r).Value = $vals[2]
    $ws.Range("E" + $r).Value = $vals[3]
    $ws.Range("F" + $r).Value = $vals[4]
    $ws.Range("G" + $r).Value = $vals[5]
    $ws.Range("H" + $r).Value = $vals[6]
}

# Re-sort the whole data table by column B ("Casos totales") descending,
# exactly like the published ranking table.
$dataRange = $ws.Range("A" + $firstDataRow + ":H" + $lastDataRow)
$sortKey = $ws.Range("B" + $firstDataRow)
$dataRange.Sort($sortKey, 2)

# Refresh the "last updated" banner in A1.
$ws.Range("A1").Value = "Datos actualizados a 6 de Agosto de 2020 a las 11:02"

Write-Output "Update complete"
